# Protokoll von 20160329 in Arbeit verbessern
#
# This script applies, via Word COM-interop, the same edits that were made
# to the "in Bearbeitung" (work-in-progress) meeting protocol:
#   - moves the hidden "_GoBack" bookmark from the old edit location to the
#     spot where the author's last edit actually happened (end of the
#     "Sonstiges" entry, row 6, of the agenda table)
#   - relabels / shifts the agenda items in the agenda table (table 2)
#   - clears out the now-unused last agenda row
#   - refreshes the cached footer fields (last-saved date/time + page number)

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (its old location, right after the
#    first table, becomes a plain empty paragraph again).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
}

# ---------------------------------------------------------------------
# 2) Update the agenda table (2nd table in the document) rows 3..8
#    (table row 2 is "Nr./Beschreibung" header, data rows start at 2).
# ---------------------------------------------------------------------
$agenda = $d.Tables.Item(2)

# Row 3 (item "2"): "Status der Aufgaben" -> "Status des Pflichtenhefts"
# (the prefix is split into two runs, "Status" + " des ", to mirror the
# author's edit).
$p = $agenda.Cell(3, 2).Range.Paragraphs.Item(1)
$xml = '<w:p ' + $wNs + '><w:r><w:t>Status</w:t></w:r><w:r><w:t xml:space="preserve"> des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pflichtenhefts</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$p.Range.InsertXML($xml)

# Row 4 (item "3"): "Status des Pflichtenhefts" -> "Status des Arbeitspakets"
$p = $agenda.Cell(4, 2).Range.Paragraphs.Item(1)
$xml = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Status des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Arbeitspakets</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$p.Range.InsertXML($xml)

# Row 5 (item "4"): "Status der Arbeitspakete" -> "Status der Aufgabe"
$p = $agenda.Cell(5, 2).Range.Paragraphs.Item(1)
$xml = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Status der </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Aufgabe</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$p.Range.InsertXML($xml)

# Row 6 (item "5"): "Werkzeugwechsel" -> "Aufgabeverteilung"
$p = $agenda.Cell(6, 2).Range.Paragraphs.Item(1)
$xml = '<w:p ' + $wNs + '><w:proofErr w:type="spellStart"/><w:r><w:t>Aufgabeverteilung</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$p.Range.InsertXML($xml)

# Row 7 (item "6"): "Feedback" -> "Sonstiges", and this is where the new
# "_GoBack" bookmark (last edit point) now lives.
$p = $agenda.Cell(7, 2).Range.Paragraphs.Item(1)
$xml = '<w:p ' + $wNs + '><w:r><w:t>Sonstiges</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$p.Range.InsertXML($xml)

# Row 8 (item "7"): old "Sonstiges" entry is cleared out entirely.
$p = $agenda.Cell(8, 2).Range.Paragraphs.Item(1)
$xml = '<w:p ' + $wNs + '/>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) Refresh the footer's cached field results: the "last saved" date and
#    time, and the current page number.
# ---------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$fp = $footer.Range.Paragraphs.Item(1)
$footerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0A2520F6" w14:textId="77777777" w:rsidR="000335B3" w:rsidRDefault="007A4D75"><w:pPr><w:pStyle w:val="HeaderFooter"/><w:tabs><w:tab w:val="clear" w:pos="9360"/><w:tab w:val="center" w:pos="4816"/><w:tab w:val="right" w:pos="9632"/></w:tabs></w:pPr><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">Version vom </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:instrText xml:space="preserve"> DATE \@ "d. MMMM y HH:mm" </w:instrText></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00F76063"><w:rPr><w:noProof/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>29. M</w:t></w:r><w:r w:rsidR="00F76063"><w:rPr><w:noProof/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>&#228;</w:t></w:r><w:r w:rsidR="00F76063"><w:rPr><w:noProof/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>rz 16 19:37</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:tab/><w:t xml:space="preserve">Seite </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:instrText xml:space="preserve"> PAGE </w:instrText></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00F76063"><w:rPr><w:noProof/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> von </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:instrText xml:space="preserve"> NUMPAGES </w:instrText></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00F76063"><w:rPr><w:noProof/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'
$fp.Range.InsertXML($footerXml)

Write-Output "done"
